$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 20
$ws.Cells.Item(2, 2).Value = -20.60755986554534
$ws.Cells.Item(2, 3).Value = 1.089226581658092
$ws.Cells.Item(2, 4).Value = 75.00302793421463
$ws.Cells.Item(2, 5).Value = -0.406779944259445
$ws.Cells.Item(2, 6).Value = -0.352054911594391
$ws.Cells.Item(2, 7).Value = 0.7815408847023604
$ws.Cells.Item(2, 8).Value = -0.4500861918092851
$ws.Cells.Item(2, 9).Value = 0.8988809271980212
$ws.Cells.Item(2, 10).Value = 4.3552525155391
$ws.Cells.Item(2, 11).Value = 27
$ws.Cells.Item(2, 12).Value = -13.08276063154555
$ws.Cells.Item(2, 13).Value = 0
$ws.Cells.Item(2, 14).Value = 4.35525251498856
$ws.Cells.Item(2, 15).Value = 5.021919181655227

$ws.Cells.Item(3, 1).Value = 21
$ws.Cells.Item(3, 2).Value = [double]"6.489415692342842e-06"
$ws.Cells.Item(3, 3).Value = 0.7035650081437954
$ws.Cells.Item(3, 4).Value = 1.195271404766907
$ws.Cells.Item(3, 5).Value = 62.69770586818203
$ws.Cells.Item(3, 6).Value = 1.653471968311442
$ws.Cells.Item(3, 7).Value = -1.727242461161577
$ws.Cells.Item(3, 8).Value = 0.5345901828952315
$ws.Cells.Item(3, 9).Value = -0.6299647063015685
$ws.Cells.Item(3, 10).Value = 4.355350724842692
$ws.Cells.Item(3, 11).Value = 98
$ws.Cells.Item(3, 12).Value = -10.91721351340661
$ws.Cells.Item(3, 13).Value = 0
$ws.Cells.Item(3, 14).Value = 4.355350721035633
$ws.Cells.Item(3, 15).Value = 5.0220173877023

$ws.Cells.Item(4, 1).Value = 22
$ws.Cells.Item(4, 2).Value = -1.836450284879223
$ws.Cells.Item(4, 3).Value = 62.8767700094226
$ws.Cells.Item(4, 4).Value = 0.1799469150252427
$ws.Cells.Item(4, 5).Value = -0.0001415292567473507
$ws.Cells.Item(4, 6).Value = -0.6433970046019519
$ws.Cells.Item(4, 7).Value = -0.7616181374878477
$ws.Cells.Item(4, 8).Value = 0.8604019339048006
$ws.Cells.Item(4, 9).Value = 1.848922505396178
$ws.Cells.Item(4, 10).Value = 4.355362752407631
$ws.Cells.Item(4, 11).Value = 17
$ws.Cells.Item(4, 12).Value = -3.972176334657421
$ws.Cells.Item(4, 13).Value = 0
$ws.Cells.Item(4, 14).Value = 4.35536275239051
$ws.Cells.Item(4, 15).Value = 5.022029419057177

$ws.Cells.Item(5, 1).Value = 23
$ws.Cells.Item(5, 2).Value = 0.3725194491593553
$ws.Cells.Item(5, 3).Value = 8.517556335889296
$ws.Cells.Item(5, 4).Value = 58.72967226374491
$ws.Cells.Item(5, 5).Value = -0.0290272403757472
$ws.Cells.Item(5, 6).Value = 0.7927274845000563
$ws.Cells.Item(5, 7).Value = -1.324838884042288
$ws.Cells.Item(5, 8).Value = -0.7077263909663822
$ws.Cells.Item(5, 9).Value = 1.073934080018848
$ws.Cells.Item(5, 10).Value = 4.355367582331654
$ws.Cells.Item(5, 11).Value = 77
$ws.Cells.Item(5, 12).Value = -5.976831165322857
$ws.Cells.Item(5, 13).Value = 0
$ws.Cells.Item(5, 14).Value = 4.355367582335208
$ws.Cells.Item(5, 15).Value = 5.022034249001875

$ws.Cells.Item(6, 1).Value = 16
$ws.Cells.Item(6, 2).Value = -0.3630679222453754
$ws.Cells.Item(6, 3).Value = 1.444344719566855
$ws.Cells.Item(6, 4).Value = 47.84348870804013
$ws.Cells.Item(6, 5).Value = 20.08156574257011
$ws.Cells.Item(6, 6).Value = 0.8703964381228366
$ws.Cells.Item(6, 7).Value = 0.7040910487935337
$ws.Cells.Item(6, 8).Value = -0.3947333488908169
$ws.Cells.Item(6, 9).Value = -0.8932278747492548
$ws.Cells.Item(6, 10).Value = 4.354054096912488
$ws.Cells.Item(6, 11).Value = 40
$ws.Cells.Item(6, 12).Value = -18.49099562563747
$ws.Cells.Item(6, 13).Value = 0
$ws.Cells.Item(6, 14).Value = 4.355476733141792
$ws.Cells.Item(6, 15).Value = 5.022143399808459

$ws.Cells.Item(7, 1).Value = 19
$ws.Cells.Item(7, 2).Value = 91.69643404025862
$ws.Cells.Item(7, 3).Value = -0.5686553812968771
$ws.Cells.Item(7, 4).Value = 1.817496675291703
$ws.Cells.Item(7, 5).Value = 558.5564966879476
$ws.Cells.Item(7, 6).Value = -1.469611940034627
$ws.Cells.Item(7, 7).Value = 0.9017877235241905
$ws.Cells.Item(7, 8).Value = 0.7566475795356018
$ws.Cells.Item(7, 9).Value = -0.01278018778117751
$ws.Cells.Item(7, 10).Value = 4.354218165535372
$ws.Cells.Item(7, 11).Value = 74
$ws.Cells.Item(7, 12).Value = -543.2311218206245
$ws.Cells.Item(7, 13).Value = 0
$ws.Cells.Item(7, 14).Value = 4.355514536895612
$ws.Cells.Item(7, 15).Value = 5.022181203562279

$ws.Cells.Item(8, 1).Value = 24
$ws.Cells.Item(8, 2).Value = -2.22303467957997
$ws.Cells.Item(8, 3).Value = 0.002041387304143306
$ws.Cells.Item(8, 4).Value = 2.123366612837303
$ws.Cells.Item(8, 5).Value = 72.80241291235619
$ws.Cells.Item(8, 6).Value = -0.8754408538818552
$ws.Cells.Item(8, 7).Value = 1.360473970173114
$ws.Cells.Item(8, 8).Value = 0.4111161267366548
$ws.Cells.Item(8, 9).Value = -0.7137636663214746
$ws.Cells.Item(8, 10).Value = 4.355515010836472
$ws.Cells.Item(8, 11).Value = 8
$ws.Cells.Item(8, 12).Value = -11.30975084725141
$ws.Cells.Item(8, 13).Value = 0
$ws.Cells.Item(8, 14).Value = 4.35551500608134
$ws.Cells.Item(8, 15).Value = 5.022181672748007

$ws.Cells.Item(9, 1).Value = 25
$ws.Cells.Item(9, 2).Value = -14.04075639617688
$ws.Cells.Item(9, 3).Value = [double]"-4.785324361805287e-05"
$ws.Cells.Item(9, 4).Value = 8.052089433750204
$ws.Cells.Item(9, 5).Value = 58.41453540999616
$ws.Cells.Item(9, 6).Value = 0.3890331918063996
$ws.Cells.Item(9, 7).Value = 1.980658151673694
$ws.Cells.Item(9, 8).Value = 0.4982875334585
$ws.Cells.Item(9, 9).Value = -0.993170253252289
$ws.Cells.Item(9, 10).Value = 4.355557487915684
$ws.Cells.Item(9, 11).Value = 35
$ws.Cells.Item(9, 12).Value = 10.75839013608952
$ws.Cells.Item(9, 13).Value = 0
$ws.Cells.Item(9, 14).Value = 4.355557487876269
$ws.Cells.Item(9, 15).Value = 5.022224154542936

$ws.Cells.Item(10, 1).Value = 15
$ws.Cells.Item(10, 2).Value = 118.8566310237852
$ws.Cells.Item(10, 3).Value = 384.9131034325741
$ws.Cells.Item(10, 4).Value = 17.40569882423881
$ws.Cells.Item(10, 5).Value = -402.9847387035821
$ws.Cells.Item(10, 6).Value = -0.5307320830933719
$ws.Cells.Item(10, 7).Value = -1.417852201440257
$ws.Cells.Item(10, 8).Value = 0.2162672493667328
$ws.Cells.Item(10, 9).Value = -1.306742012696025
$ws.Cells.Item(10, 10).Value = 4.353988233153366
$ws.Cells.Item(10, 11).Value = 26
$ws.Cells.Item(10, 12).Value = -50.44277963598873
$ws.Cells.Item(10, 13).Value = 0
$ws.Cells.Item(10, 14).Value = 4.355620433024086
$ws.Cells.Item(10, 15).Value = 5.022287099690753

$ws.Cells.Item(11, 1).Value = 26
$ws.Cells.Item(11, 2).Value = -2.402943353734941
$ws.Cells.Item(11, 3).Value = 44.53439659643279
$ws.Cells.Item(11, 4).Value = -0.01328063671396947
$ws.Cells.Item(11, 5).Value = 2.399261227876057
$ws.Cells.Item(11, 6).Value = 0.7065799905432026
$ws.Cells.Item(11, 7).Value = -0.4627780455198596
$ws.Cells.Item(11, 8).Value = 1.301177413002297
$ws.Cells.Item(11, 9).Value = 0.7558072736687058
$ws.Cells.Item(11, 10).Value = 4.355632373509996
$ws.Cells.Item(11, 11).Value = 18
$ws.Cells.Item(11, 12).Value = -9.078821358084308
$ws.Cells.Item(11, 13).Value = 0
$ws.Cells.Item(11, 14).Value = 4.355632374199807
$ws.Cells.Item(11, 15).Value = 5.022299040866474
